{"js": "// Word JS API (Office.js) edit script.\n// Body of: async (context) => { ... }\n//\n// Summary of the edit (per the XML diff):\n//  1. Opening line \"I hope this email finds you well.\" becomes a longer,\n//     reworded greeting/introduction.\n//  2. The paragraph that used to start \"I was previously an undergraduate\n//     student ...\" is replaced with a new sentence about having learned of\n//     the recipient's career, while keeping the trailing\n//     \"Currently, I have several papers under review ... submission to PNAS.\"\n//     sentence intact.\n//  3. \"to high-dimensional.\" gets extended with \", which is the fourth\n//     paper of REDS series.\"\n//  4. The trailing empty paragraph after \"Tuobang Li\" is removed.\n\nconst body = context.document.body;\n\n// --- 1. Reword the opening greeting -----------------------------------\nconst greetingResults = body.search(\"I hope this email finds you well.\", { matchCase: true });\ngreetingResults.load(\"items\");\nawait context.sync();\n\ngreetingResults.items[0].insertText(\n  \"I hope this message finds you well. I am writing to you as a former student of Guangdong Technion who has recently transferred to Berkeley due to personal reasons.\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- 2. Replace the \"I was previously ...\" sentence with the new one ---\nconst introResults = body.search(\n  \"I was previously an undergraduate student from Guangdong Technion. I am going to arrive Berkeley soon. Currently,\",\n  { matchCase: true }\n);\nintroResults.load(\"items\");\nawait context.sync();\n\nintroResults.items[0].insertText(\n  \"I have learned of your esteemed career and your connection to Technion. Currently,\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- 3. Extend the REDS sentence ---------------------------------------\nconst redsResults = body.search(\"to high-dimensional. The basic principle\", { matchCase: true });\nredsResults.load(\"items\");\nawait context.sync();\n\nredsResults.items[0].insertText(\n  \"to high-dimensional, which is the fourth paper of REDS series. The basic principle\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- 4. Remove the trailing empty paragraph after \"Tuobang Li\" ---------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst count = paragraphs.items.length;\nconst lastParagraph = paragraphs.items[count - 1];\nconst secondToLastParagraph = paragraphs.items[count - 2];\n\n// Collapse the break between the \"Tuobang Li\" paragraph and the empty\n// paragraph that follows it by deleting the range spanning from the end\n// of the former to the end of the latter (the very last paragraph mark of\n// the body cannot itself be deleted, so we must merge into it instead of\n// calling .delete() on it directly).\nconst mergeStart = secondToLastParagraph.getRange(\"End\");\nconst mergeEnd = lastParagraph.getRange(\"End\");\nconst mergeRange = mergeStart.expandTo(mergeEnd);\nmergeRange.delete();\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word / $doc / $d (ActiveDocument) are pre-seeded by the host.\n#\n# Summary of the edit (per the XML diff):\n#  1. Opening line \"I hope this email finds you well.\" becomes a longer,\n#     reworded greeting/introduction.\n#  2. The paragraph that used to start \"I was previously an undergraduate\n#     student ...\" is replaced with a new sentence about having learned of\n#     the recipient's career, while keeping the trailing\n#     \"Currently, I have several papers under review ... submission to PNAS.\"\n#     sentence intact.\n#  3. \"to high-dimensional.\" gets extended with \", which is the fourth\n#     paper of REDS series.\"\n#  4. The trailing empty paragraph after \"Tuobang Li\" is removed.\n\n$d = $word.ActiveDocument\n\n# --- 1. Reword the opening greeting -----------------------------------\n$range1 = $d.Content\n$range1.Find.Execute(\n    \"I hope this email finds you well.\",  # FindText\n    $true,                                 # MatchCase\n    $true,                                 # MatchWholeWord\n    $false,                                # MatchWildcards\n    $false,                                # MatchSoundsLike\n    $false,                                # MatchAllWordForms\n    $true,                                 # Forward\n    1,                                      # Wrap (wdFindContinue)\n    $false,                                # Format\n    \"I hope this message finds you well. I am writing to you as a former student of Guangdong Technion who has recently transferred to Berkeley due to personal reasons.\",  # ReplaceWith\n    2                                       # Replace (wdReplaceOne)\n) | Out-Null\n\n# --- 2. Replace the \"I was previously ...\" sentence with the new one ---\n$range2 = $d.Content\n$range2.Find.Execute(\n    \"I was previously an undergraduate student from Guangdong Technion. I am going to arrive Berkeley soon. Currently,\",\n    $true,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"I have learned of your esteemed career and your connection to Technion. Currently,\",\n    2\n) | Out-Null\n\n# --- 3. Extend the REDS sentence ---------------------------------------\n$range3 = $d.Content\n$range3.Find.Execute(\n    \"to high-dimensional. The basic principle\",\n    $true,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"to high-dimensional, which is the fourth paper of REDS series. The basic principle\",\n    2\n) | Out-Null\n\n# --- 4. Remove the trailing empty paragraph after \"Tuobang Li\" ---------\n$count = $d.Paragraphs.Count\n$secondToLast = $d.Paragraphs.Item($count - 1)\n\n# The very last paragraph mark of the body story cannot be deleted on its\n# own (Word keeps at least one paragraph mark), so collapse the range to\n# the end of the \"Tuobang Li\" paragraph and extend it by one character to\n# swallow the following (empty) paragraph's mark instead.\n$mergeRange = $secondToLast.Range\n$mergeRange.SetRange($mergeRange.End, $mergeRange.End)\n$mergeRange.MoveEnd(1, 1) | Out-Null\n$mergeRange.Delete()\n"}
